# Avg Vehicle Loading.xlsx - calibrate transp sector ships, rail, aviation
#
# This script:
#  1. On "AVLo-passengers": converts the formula-driven rows for aircraft (4),
#     rail (5) and ships (6) into hard static values, highlighting the
#     overridden cells with a pale-yellow fill (the classic "manual
#     calibration override" color). The rail row value is recalibrated
#     (x10).
#  2. On "AVLo-freight": converts the same three rows (aircraft, rail,
#     ships) into hard static values with the same highlight, keeping the
#     existing computed magnitudes.
#  3. Updates a couple of selections / the active sheet to match what was
#     left selected when the author saved the file.

$wb = $excel.ActiveWorkbook

# Color used to flag manually-overridden calibration cells: RGB(255,255,153)
$calibrationColor = 10092543

# ----------------------------------------------------------------------
# Sheet: AVLo-passengers
# ----------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("AVLo-passengers")

# Row 4 - aircraft: freeze formulas to values (value itself unchanged)
$rowB = $wsP.Range("B4")
$rowB.Value = $rowB.Value2
$rowB.NumberFormat = "0.0"
$rowB.Interior.Color = $calibrationColor

$rowRest = $wsP.Range("C4:AK4")
$rowRest.Value = 111.39416306433705
$rowRest.NumberFormat = "0.00"
$rowRest.Interior.Color = $calibrationColor

# Row 5 - rail: recalibrated to 10x the previous value, freeze to values
$railB = $wsP.Range("B5")
$railB.Value = 486.56731685074101
$railB.NumberFormat = "0.0"
$railB.Interior.Color = $calibrationColor

$railRest = $wsP.Range("C5:AK5")
$railRest.Value = 486.56731685074101
$railRest.NumberFormat = "0.00"
$railRest.Interior.Color = $calibrationColor

# Row 6 - ships: freeze formulas to values (value itself unchanged)
$shipB = $wsP.Range("B6")
$shipB.Value = 1
$shipB.NumberFormat = "0.00"
$shipB.Interior.Color = $calibrationColor

$shipRest = $wsP.Range("C6:AK6")
$shipRest.Value = 1
$shipRest.NumberFormat = "0.00"
$shipRest.Interior.Color = $calibrationColor

# Selection left on this sheet after editing
$wsP.Range("C15").Select()

# ----------------------------------------------------------------------
# Sheet: AVLo-freight
# ----------------------------------------------------------------------
$wsF = $wb.Worksheets.Item("AVLo-freight")

# Row 4 - aircraft
$fAirB = $wsF.Range("B4")
$fAirB.Value = $fAirB.Value2
$fAirB.NumberFormat = "0"
$fAirB.Interior.Color = $calibrationColor

$fAirRest = $wsF.Range("C4:AJ4")
$fAirRest.Value = 41.989116133258747
$fAirRest.NumberFormat = "0"
$fAirRest.Interior.Color = $calibrationColor

# Row 5 - rail
$fRailB = $wsF.Range("B5")
$fRailB.Value = 3512.35916421195
$fRailB.NumberFormat = "0"
$fRailB.Interior.Color = $calibrationColor

$fRailRest = $wsF.Range("C5:AJ5")
$fRailRest.Value = 3512.35916421195
$fRailRest.NumberFormat = "0"
$fRailRest.Interior.Color = $calibrationColor

# Row 6 - ships
$fShipB = $wsF.Range("B6")
$fShipB.Value = 1974.4736422180429
$fShipB.NumberFormat = "0"
$fShipB.Interior.Color = $calibrationColor

$fShipRest = $wsF.Range("C6:AJ6")
$fShipRest.Value = 1974.4736422180429
$fShipRest.NumberFormat = "0"
$fShipRest.Interior.Color = $calibrationColor

# Selection left on this sheet, and this sheet is the one left active/selected
$wsF.Range("D19").Select()
$wsF.Activate()

# ----------------------------------------------------------------------
# Minor leftover selections on other sheets (cosmetic, matches saved state)
# ----------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B3").Select()

$wsFreightAct = $wb.Worksheets.Item("Mexico Freight Activity")
$wsFreightAct.Range("U19").Select()

$wsLDV = $wb.Worksheets.Item("Mexico Psgr LDVs, Psgr HDVs")
$wsLDV.Range("B16").Select()

# Re-activate AVLo-freight last so it ends up the selected/visible tab
$wsF.Activate()
$wsF.Range("D19").Select()
